$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Efna2", "Epha3", "ECs", 3, 1, 3.099562666666667, 9.298688, 0.2686390288432488, 0.2686390288432488, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.2349344518826667, 2.114410066944, 0.0004782314557535506, 0.0004782314557535506),
    @("ECs", "Efna2", "Epha3", "FAPs", 3, 1, 3.099562666666667, 9.298688, 0.2686390288432488, 0.2686390288432488, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 124.3772060098275, 1119.394854088448, 0.253181650524146, 0.2531816505241459),
    @("ECs", "Efna2", "Epha3", "MuSCs", 3, 1, 3.099562666666667, 9.298688, 0.2686390288432488, 0.2686390288432488, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 7.265600125553778, 65.390401129984, 0.0147898211485053, 0.0147898211485053),
    @("ECs", "Efna2", "Epha3", "Resolving-Mac", 3, 1, 3.099562666666667, 9.298688, 0.2686390288432488, 0.2686390288432488, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.09300754375111112, 0.83706789376, 0.0001893257148439435, 0.0001893257148439435),
    @("FAPs", "Efna2", "Epha3", "ECs", 3, 1, 6.189892666666666, 18.569678, 0.5364778626674904, 0.5364778626674905, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.4691691045626666, 4.222521941064, 0.0009550384035699102, 0.0009550384035699104),
    @("FAPs", "Efna2", "Epha3", "FAPs", 3, 1, 6.189892666666666, 18.569678, 0.5364778626674904, 0.5364778626674905, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 248.3839296621375, 2235.455366959238, 0.5056091489188498, 0.5056091489188498),
    @("FAPs", "Efna2", "Epha3", "MuSCs", 3, 1, 6.189892666666666, 18.569678, 0.5364778626674904, 0.5364778626674905, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 14.50955820953378, 130.586023885804, 0.02953558785985007, 0.02953558785985007),
    @("FAPs", "Efna2", "Epha3", "Resolving-Mac", 3, 1, 6.189892666666666, 18.569678, 0.5364778626674904, 0.5364778626674905, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.1857380459511111, 1.67164241356, 0.0003780874852206948, 0.0003780874852206948),
    @("MuSCs", "Efna2", "Epha3", "ECs", 3, 1, 1.888584, 5.665752, 0.1636835341659699, 0.1636835341659699, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.143147112864, 1.288324015776, 0.0002913895838744768, 0.0002913895838744768),
    @("MuSCs", "Efna2", "Epha3", "FAPs", 3, 1, 1.888584, 5.665752, 0.1636835341659699, 0.1636835341659699, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 75.783852916088, 682.0546762447921, 0.1542652514871433, 0.1542652514871433),
    @("MuSCs", "Efna2", "Epha3", "MuSCs", 3, 1, 1.888584, 5.665752, 0.1636835341659699, 0.1636835341659699, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 4.426978133104, 39.842803197936, 0.009011535686731954, 0.009011535686731955),
    @("MuSCs", "Efna2", "Epha3", "Resolving-Mac", 3, 1, 1.888584, 5.665752, 0.1636835341659699, 0.1636835341659699, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.05667011056, 0.51003099504, 0.0001153574082202245, 0.0001153574082202245),
    @("Resolving-Mac", "Efna2", "Epha3", "ECs", 2, 0.6666666666666666, 0.3599813333333333, 1.079944, 0.03119957432329092, 0.03119957432329093, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.02728514514133333, 0.245566306272, 0.00005554151201248095, 0.00005554151201248096),
    @("Resolving-Mac", "Efna2", "Epha3", "FAPs", 2, 0.6666666666666666, 0.3599813333333333, 1.079944, 0.03119957432329092, 0.03119957432329093, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 14.44509347631378, 130.005841286824, 0.02940436375471984, 0.02940436375471984),
    @("Resolving-Mac", "Efna2", "Epha3", "MuSCs", 2, 0.6666666666666666, 0.3599813333333333, 1.079944, 0.03119957432329092, 0.03119957432329093, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 0.8438224039768888, 7.594401635792, 0.001717680882550463, 0.001717680882550464),
    @("Resolving-Mac", "Efna2", "Epha3", "Resolving-Mac", 2, 0.6666666666666666, 0.3599813333333333, 1.079944, 0.03119957432329092, 0.03119957432329093, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.01080183987555556, 0.09721655888, 0.00002198817400814263, 0.00002198817400814263)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $rowVals[$j]
    }
}
